$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 3290.886797766764
$ws.Range("C2").Value = 1792.209789670246
$ws.Range("D2").Value = 1661.023907426512
